$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 already contain the "1007609217" shared string with the default
# (unstyled) cell format. Copy them down into rows 6-8 so the new cells
# reuse the same shared-string entry and keep the plain/default cell style
# (typing the numeric-looking text directly would make Excel store it as a
# number instead of reusing the shared string).
$ws.Range("A2:A4").Copy()
$ws.Range("A6:A8").PasteSpecial(-4104)
